$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($addr, $value) {
    $ws.Range($addr).Value = $value
}

function Set-ForcedText($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-PlainText "D2" "38.470.25"
Set-PlainText "E2" "  +2.34%  "

# Row 3 - Ethereum
Set-PlainText "D3" "2.081.12"
Set-PlainText "E3" "  +2.52%  "

# Row 4 - TetherUSD
Set-PlainText "E4" "  -0.03%  "

# Row 5 - BNB
Set-ForcedText "D5" "228.42"
Set-PlainText "E5" "  +0.96%  "

# Row 6 - XRP
Set-PlainText "E6" "  +0.76%  "

# Row 7 - Solana
Set-ForcedText "D7" "60.35"
Set-PlainText "E7" "  +1.35%  "

# Row 8 - USDC
Set-PlainText "E8" "  -0.01%  "

# Row 9 - Cardano
Set-ForcedText "D9" "0.382"

# Row 10 - Dogecoin
Set-ForcedText "D10" "0.0834"
Set-PlainText "E10" "  +1.12%  "

# Row 11 - TRON
Set-PlainText "E11" "  -0.39%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-PlainText "D12" "2.389.78"
Set-PlainText "E12" "  +2.63%  "

# Row 13 - Chainlink
Set-PlainText "E13" "  +4.21%  "

# Row 14 - Avalanche
Set-ForcedText "D14" "22.28"
Set-PlainText "E14" "  +6.02%  "

# Row 15 - Polygon
Set-ForcedText "D15" "0.783"
Set-PlainText "E15" "  +1.72%  "

# Row 16 - Polkadot
Set-PlainText "E16" "  -0.87%  "

# Row 17 - WrappedEther
Set-PlainText "D17" "2.077.50"
Set-PlainText "E17" "  +1.94%  "

# Row 18 - WrappedBTC
Set-PlainText "D18" "38.402.79"
Set-PlainText "E18" "  +2.23%  "

# Row 19 - Litecoin
Set-ForcedText "D19" "71.41"
Set-PlainText "E19" "  +3.15%  "

# Row 20 - Uniswap
Set-PlainText "E20" "  +0.96%  "

# Row 22 - BitcoinCash
Set-ForcedText "D22" "225.04"
Set-PlainText "E22" "  +0.62%  "

# Row 24 - Toncoin
Set-PlainText "E24" "  +0.10%  "

# Row 25 - PancakeSwap
Set-ForcedText "D25" "2.32"
Set-PlainText "E25" "  +2.70%  "

# Row 26 - Monero
Set-ForcedText "D26" "170.68"
Set-PlainText "E26" "  +1.88%  "

# Row 27 - Cosmos
Set-PlainText "E27" "  +0.68%  "

# Row 28 - Kaspa
Set-ForcedText "D28" "0.135"
Set-PlainText "E28" "  +6.42%  "

# Row 29 - ImmutableX
Set-PlainText "E29" "  +9.17%  "

# Row 30 - EthereumClassic
Set-ForcedText "D30" "19.08"
Set-PlainText "E30" "  +1.95%  "

# Row 31 - Stellar
Set-PlainText "E31" "  -0.01%  "

# Row 32 - WEMIXToken
Set-PlainText "E32" "  +4.65%  "

# Row 33 - InternetComputer(DFINITY)
Set-PlainText "E33" "  +7.48%  "

# Row 34 - Filecoin
Set-ForcedText "D34" "4.48"
Set-PlainText "E34" "  +2.65%  "

# Row 35 - Hedera
Set-PlainText "E35" "  +0.52%  "

# Row 36 - LidoDAOToken
Set-PlainText "E36" "  +1.22%  "

# Row 37 - THORChain
Set-ForcedText "D37" "6.32"
Set-PlainText "E37" "  -2.38%  "

# Row 38 - RenderToken
Set-PlainText "E38" "  +4.12%  "

# Row 39 - BinanceUSD
Set-ForcedText "D39" "1.00"
Set-PlainText "E39" "  +0.00%  "

# Row 40 - InjectiveProtocol
Set-ForcedText "D40" "18.29"
Set-PlainText "E40" "  +2.30%  "

# Row 41 - Maker
Set-PlainText "D41" "1.540.65"
Set-PlainText "E41" "  +1.28%  "

# Row 42 - Aave
Set-ForcedText "D42" "100.38"
Set-PlainText "E42" "  +3.93%  "

# Row 43 - VeChain
Set-PlainText "E43" "  +3.20%  "

# Row 44 - HuobiToken -> Cronos
Set-PlainText "B44" "Cronos"
Set-PlainText "C44" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-ForcedText "D44" "0.0922"
Set-PlainText "E44" "  +1.85%  "

# Row 45 - Cronos -> HuobiToken
Set-PlainText "B45" "HuobiToken"
Set-PlainText "C45" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-ForcedText "D45" "2.81"
Set-PlainText "E45" "  -0.89%  "

# Row 46 - FraxShare
Set-ForcedText "D46" "7.67"
Set-PlainText "E46" "  +8.86%  "

# Row 47 - FTXToken
Set-ForcedText "D47" "4.11"
Set-PlainText "E47" "  -2.80%  "

# Row 48 - TrustWalletToken
Set-PlainText "E48" "  +1.12%  "

# Row 49 - ARBITRUM
Set-PlainText "E49" "  +2.56%  "

# Row 50 - MXToken
Set-ForcedText "D50" "2.99"
Set-PlainText "E50" "  +1.86%  "

# Row 51 - RocketPoolETH
Set-PlainText "D51" "2.278.18"
Set-PlainText "E51" "  +2.57%  "
